$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("C6").Value = 14174.5
$ws.Range("D6").Value = 1118.9100000000001
$ws.Range("E6").Value = 394.29599999999999
$ws.Range("F6").Value = 3721
$ws.Range("G6").Value = 8000
$ws.Range("I6").Value = 676

$ws.Range("I7").Select()
